# Generate Report for Handback
# Adds "Latest Target File" (E) / "Latest Handback File" (F) hyperlink cells
# for the two localized source files on the zh-cn and de-de sheets, updates
# the handback status text, and stamps the new "Latest Handback DateTime".

$wb = $excel.ActiveWorkbook

$sheetsInfo = @(
    @{ Name = "zh-cn";
       MdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/ec2c27aa40c768ce58c34fa387b0a692a1e403ed/e2e/7adb2795-059c-4e3f-89a2-e50de085827b.md";
       MdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/ec2c27aa40c768ce58c34fa387b0a692a1e403ed/e2e/9ef8f679-eb76-4707-8fa4-8517cddb2476.md";
       XlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/239e8a44d852c575c06db582199ef09bec5347e5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/7adb2795-059c-4e3f-89a2-e50de085827b.0c77cf34e51e3f0ff36c8300b42ffca9402a27d7.zh-cn.xlf";
       XlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/239e8a44d852c575c06db582199ef09bec5347e5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9ef8f679-eb76-4707-8fa4-8517cddb2476.75243989c2ad2794e108b88e6ec4b39f93d5be2f.zh-cn.xlf";
       MdName = "7adb2795-059c-4e3f-89a2-e50de085827b.md";
       MdName2 = "9ef8f679-eb76-4707-8fa4-8517cddb2476.md";
       XlfName = "7adb2795-059c-4e3f-89a2-e50de085827b.0c77cf34e51e3f0ff36c8300b42ffca9402a27d7.zh-cn.xlf";
       XlfName2 = "9ef8f679-eb76-4707-8fa4-8517cddb2476.75243989c2ad2794e108b88e6ec4b39f93d5be2f.zh-cn.xlf";
       HandbackTime = "2016-03-10 22:57:58" },
    @{ Name = "de-de";
       MdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/ec2c27aa40c768ce58c34fa387b0a692a1e403ed/e2e/7adb2795-059c-4e3f-89a2-e50de085827b.md";
       MdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/ec2c27aa40c768ce58c34fa387b0a692a1e403ed/e2e/9ef8f679-eb76-4707-8fa4-8517cddb2476.md";
       XlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6f254743196fa4072732762e8e930a8b1b487db7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/7adb2795-059c-4e3f-89a2-e50de085827b.0c77cf34e51e3f0ff36c8300b42ffca9402a27d7.de-de.xlf";
       XlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6f254743196fa4072732762e8e930a8b1b487db7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9ef8f679-eb76-4707-8fa4-8517cddb2476.75243989c2ad2794e108b88e6ec4b39f93d5be2f.de-de.xlf";
       MdName = "7adb2795-059c-4e3f-89a2-e50de085827b.md";
       MdName2 = "9ef8f679-eb76-4707-8fa4-8517cddb2476.md";
       XlfName = "7adb2795-059c-4e3f-89a2-e50de085827b.0c77cf34e51e3f0ff36c8300b42ffca9402a27d7.de-de.xlf";
       XlfName2 = "9ef8f679-eb76-4707-8fa4-8517cddb2476.75243989c2ad2794e108b88e6ec4b39f93d5be2f.de-de.xlf";
       HandbackTime = "2016-03-10 22:58:18" }
)

foreach ($info in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Status for both handed-off rows is now "Handed back: in sync with en-US"
    $ws.Range("B2").Value = "Handed back: in sync with en-US"
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    # Row 2 - Latest Target File / Latest Handback File
    $ws.Range("E2").Value = $info.MdName
    $ws.Hyperlinks.Add($ws.Range("E2"), $info.MdUrl, "", "", $info.MdName) | Out-Null

    $ws.Range("F2").Value = $info.XlfName
    $ws.Hyperlinks.Add($ws.Range("F2"), $info.XlfUrl, "", "", $info.XlfName) | Out-Null

    # Row 3 - Latest Target File / Latest Handback File
    $ws.Range("E3").Value = $info.MdName2
    $ws.Hyperlinks.Add($ws.Range("E3"), $info.MdUrl2, "", "", $info.MdName2) | Out-Null

    $ws.Range("F3").Value = $info.XlfName2
    $ws.Hyperlinks.Add($ws.Range("F3"), $info.XlfUrl2, "", "", $info.XlfName2) | Out-Null

    # Latest Handback DateTime for both rows
    $ws.Range("G2").Value = $info.HandbackTime
    $ws.Range("G3").Value = $info.HandbackTime
}
